$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (stimRamp, rampLin, logStep) right after gabNum (col I),
# pushing the old J..S columns to M..V.
$ws.Columns("J:L").Insert()

# Insert one new column (maskRamp) right after the (now shifted) maskRR column (col N),
# pushing the old O..S (now at that position) further right.
$ws.Columns("O").Insert()

# Header row (row 1) labels for the newly inserted columns.
$ws.Range("J1").Value = "stimRamp"
$ws.Range("K1").Value = "rampLin"
$ws.Range("L1").Value = "logStep"
$ws.Range("O1").Value = "maskRamp"

# Data rows 2-6: fill in values for the new columns, and bump stimT (col F) 1000 -> 2000.
for ($r = 2; $r -le 6; $r++) {
    $ws.Range("F$r").Value = 2000
    $ws.Range("J$r").Value = 1
    $ws.Range("K$r").Value = 1
    $ws.Range("L$r").Value = 1
    $ws.Range("O$r").Value = 0
}

# Selection moves to F2:F6 (matches the edited stimT column).
$ws.Range("F2:F6").Select() | Out-Null
